# "M14 Froze Encoder 12" - refresh the per-epoch accuracy values in column B
# (B3:B118) on the "Epoch Accuracy" sheet to reflect the latest training run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row (sheet row number) -> new accuracy value for column B.
# Rows not listed here (8, 14, 104, 110, 112) are unchanged from the prior run.
$values = @{
    3 = 0.421875
    4 = 0.40625
    5 = 0.359375
    6 = 0.34375
    7 = 0.28125
    9 = 0.28125
    10 = 0.265625
    11 = 0.265625
    12 = 0.3125
    13 = 0.296875
    15 = 0.34375
    16 = 0.28125
    17 = 0.28125
    18 = 0.234375
    19 = 0.203125
    20 = 0.203125
    21 = 0.140625
    22 = 0.15625
    23 = 0.203125
    24 = 0.171875
    25 = 0.15625
    26 = 0.171875
    27 = 0.171875
    28 = 0.140625
    29 = 0.171875
    30 = 0.140625
    31 = 0.15625
    32 = 0.140625
    33 = 0.171875
    34 = 0.171875
    35 = 0.171875
    36 = 0.15625
    37 = 0.15625
    38 = 0.15625
    39 = 0.15625
    40 = 0.15625
    41 = 0.15625
    42 = 0.15625
    43 = 0.15625
    44 = 0.15625
    45 = 0.15625
    46 = 0.15625
    47 = 0.15625
    48 = 0.15625
    49 = 0.15625
    50 = 0.15625
    51 = 0.15625
    52 = 0.15625
    53 = 0.15625
    54 = 0.15625
    55 = 0.15625
    56 = 0.15625
    57 = 0.15625
    58 = 0.15625
    59 = 0.15625
    60 = 0.15625
    61 = 0.15625
    62 = 0.15625
    63 = 0.15625
    64 = 0.15625
    65 = 0.15625
    66 = 0.15625
    67 = 0.15625
    68 = 0.15625
    69 = 0.15625
    70 = 0.15625
    71 = 0.15625
    72 = 0.15625
    73 = 0.15625
    74 = 0.15625
    75 = 0.15625
    76 = 0.15625
    77 = 0.15625
    78 = 0.15625
    79 = 0.15625
    80 = 0.15625
    81 = 0.15625
    82 = 0.15625
    83 = 0.15625
    84 = 0.15625
    85 = 0.15625
    86 = 0.15625
    87 = 0.15625
    88 = 0.15625
    89 = 0.15625
    90 = 0.15625
    91 = 0.15625
    92 = 0.140625
    93 = 0.140625
    94 = 0.140625
    95 = 0.140625
    96 = 0.140625
    97 = 0.140625
    98 = 0.140625
    99 = 0.140625
    100 = 0.140625
    101 = 0.140625
    102 = 0.140625
    103 = 0.125
    105 = 0.265625
    106 = 0.171875
    107 = 0.359375
    108 = 0.265625
    109 = 0.1875
    111 = 0.21875
    113 = 0.234375
    114 = 0.203125
    115 = 0.25
    116 = 0.296875
    117 = 0.28125
    118 = 0.3114754098360656
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

Write-Output "Updated $($values.Count) cells in column B"
